# Update slide deck 3
# - Refresh the "last modified" date/time fields shown on the notes pages,
#   notes master and handout master.
# - Slide 8 "Reading further": repoint the SharePoint Framework reading
#   links from docs.microsoft.com to learn.microsoft.com, bump the second
#   link's font size to match the first, and let the placeholder grow to
#   fit the now-taller (wrapped) text.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date/time fields: handout master, notes master, notes pages.
# ---------------------------------------------------------------------
$oldDate = "4/25/2022 3:38 PM"
$newDate = "11/20/2022 9:22 AM"

function Update-DateField($shape) {
    if ($null -eq $shape) { return }
    if (-not $shape.HasTextFrame) { return }
    $tr = $shape.TextFrame.TextRange
    if ($tr.Text -eq $oldDate) {
        $tr.Text = $newDate
    }
}

$hm = $p.HandoutMaster
for ($i = 1; $i -le $hm.Shapes.Count; $i++) {
    Update-DateField $hm.Shapes.Item($i)
}

$nm = $p.NotesMaster
for ($i = 1; $i -le $nm.Shapes.Count; $i++) {
    Update-DateField $nm.Shapes.Item($i)
}

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    $np = $slide.NotesPage
    for ($i = 1; $i -le $np.Shapes.Count; $i++) {
        Update-DateField $np.Shapes.Item($i)
    }
}

# ---------------------------------------------------------------------
# 2) Slide 8: "Reading further" links.
# ---------------------------------------------------------------------
$slide8 = $p.Slides.Item(8)
$linkShape = $slide8.Shapes.Item(3)
$tr = $linkShape.TextFrame.TextRange

$oldUrl1 = "https://docs.microsoft.com/sharepoint/dev/spfx/sharepoint-framework-overview"
$newUrl1 = "https://learn.microsoft.com/sharepoint/dev/spfx/sharepoint-framework-overview"
$oldUrl2 = "https://docs.microsoft.com/sharepoint/dev/spfx/integrate-with-teams-introduction"
$newUrl2 = "https://learn.microsoft.com/sharepoint/dev/spfx/integrate-with-teams-introduction"

$full = $tr.Text

$idx1 = $full.IndexOf($oldUrl1)
$tr.Characters($idx1 + 1, $oldUrl1.Length).Text = $newUrl1

# Re-locate the second URL after the first replacement (its own offset is
# unaffected since it comes later in the text, but re-query to be safe).
$full = $tr.Text
$idx2 = $full.IndexOf($oldUrl2)
$tr.Characters($idx2 + 1, $oldUrl2.Length).Text = $newUrl2

$full = $tr.Text
$newLen2 = $newUrl2.Length
$run2 = $tr.Characters($idx2 + 1, $newLen2)
$run2.Font.Size = 18

$spacePos = $idx2 + 1 + $newLen2
$trailingSpace = $tr.Characters($spacePos, 1)
$trailingSpace.Font.Size = 18

# The placeholder auto-fits to its text (spAutoFit); with the larger font
# the text now wraps onto one extra line, so the shape grows taller.
$linkShape.Height = 1692771 / 12700
